# Actualizacion Datos Personales 4 nov
# Fill in the missing contact/tutor details for the student in row 7
# (CORTEZ ANTONIO CRISTIAN JAVIER), whose record was previously incomplete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correo (student email)
$ws.Range("E7").Value = "cristianantoniof2020@gmail.com"

# Tel_Movil / Tel_Fijo (phone numbers) - force text so leading formatting
# matches the rest of the sheet, then drop back to the default style so no
# residual number-format is left applied to the cell.
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2721079631"
$ws.Range("F7").Style = "Normal"

$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "2721079631"
$ws.Range("G7").Style = "Normal"

# Tutor (tutor's name)
$ws.Range("H7").Value = "FRANCISCO JAVIER CORTÉZ LEYNES"

# Correo_Tutor (tutor's email)
$ws.Range("I7").Value = "Cristianantoniof2020@gmail.com"
